# Mean Values and Box Plot Section
# Rename the worksheet and update the active selection to match the
# state captured after the author started laying out the new section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" -> "LIME"
$ws.Name = "LIME"

# Move / update the current selection to L11 (author had clicked there
# while scoping out where the new Mean Values / Box Plot section would go)
$ws.Range("L11").Select()
